$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "21+53=74"
$t.Cell(1,2).Range.Text = "56-13=43"
$t.Cell(1,3).Range.Text = "63-17=46"
$t.Cell(1,4).Range.Text = "67+11=78"
$t.Cell(1,5).Range.Text = "23+33=56"
$t.Cell(2,1).Range.Text = "96-17=79"
$t.Cell(2,2).Range.Text = "16+19=35"
$t.Cell(2,3).Range.Text = "63-38=25"
$t.Cell(2,4).Range.Text = "38-18=20"
$t.Cell(2,5).Range.Text = "32+63=95"
$t.Cell(3,1).Range.Text = "79-48=31"
$t.Cell(3,2).Range.Text = "30+29=59"
$t.Cell(3,3).Range.Text = "75-2=73"
$t.Cell(3,4).Range.Text = "9+46=55"
$t.Cell(3,5).Range.Text = "39+36=75"
$t.Cell(4,1).Range.Text = "15+39=54"
$t.Cell(4,2).Range.Text = "7+69=76"
$t.Cell(4,3).Range.Text = "15+66=81"
$t.Cell(4,4).Range.Text = "16+10=26"
$t.Cell(4,5).Range.Text = "16+59=75"
$t.Cell(5,1).Range.Text = "37+35=72"
$t.Cell(5,2).Range.Text = "72+21=93"
$t.Cell(5,3).Range.Text = "27+14=41"
$t.Cell(5,4).Range.Text = "75-31=44"
$t.Cell(5,5).Range.Text = "64-0=64"
$t.Cell(6,1).Range.Text = "76-51=25"
$t.Cell(6,2).Range.Text = "91-23=68"
$t.Cell(6,3).Range.Text = "76-14=62"
$t.Cell(6,4).Range.Text = "9+57=66"
$t.Cell(6,5).Range.Text = "50-0=50"
$t.Cell(7,1).Range.Text = "87-53=34"
$t.Cell(7,2).Range.Text = "77-45=32"
$t.Cell(7,3).Range.Text = "88-82=6"
$t.Cell(7,4).Range.Text = "13+21=34"
$t.Cell(7,5).Range.Text = "84-79=5"
$t.Cell(8,1).Range.Text = "40-29=11"
$t.Cell(8,2).Range.Text = "82-9=73"
$t.Cell(8,3).Range.Text = "42+24=66"
$t.Cell(8,4).Range.Text = "37+19=56"
$t.Cell(8,5).Range.Text = "76-2=74"
$t.Cell(9,1).Range.Text = "99-17=82"
$t.Cell(9,2).Range.Text = "83-11=72"
$t.Cell(9,3).Range.Text = "65-48=17"
$t.Cell(9,4).Range.Text = "12-11=1"
$t.Cell(9,5).Range.Text = "93-0=93"
$t.Cell(10,1).Range.Text = "94-82=12"
$t.Cell(10,2).Range.Text = "76-43=33"
$t.Cell(10,3).Range.Text = "97-28=69"
$t.Cell(10,4).Range.Text = "0+20=20"
$t.Cell(10,5).Range.Text = "77-70=7"
$t.Cell(11,1).Range.Text = "35-15=20"
$t.Cell(11,2).Range.Text = "20+76=96"
$t.Cell(11,3).Range.Text = "2+83=85"
$t.Cell(11,4).Range.Text = "0+69=69"
$t.Cell(11,5).Range.Text = "62+33=95"
$t.Cell(12,1).Range.Text = "74-5=69"
$t.Cell(12,2).Range.Text = "27+15=42"
$t.Cell(12,3).Range.Text = "23+8=31"
$t.Cell(12,4).Range.Text = "53-23=30"
$t.Cell(12,5).Range.Text = "8+42=50"
$t.Cell(13,1).Range.Text = "52-10=42"
$t.Cell(13,2).Range.Text = "82-42=40"
$t.Cell(13,3).Range.Text = "41+35=76"
$t.Cell(13,4).Range.Text = "0+81=81"
$t.Cell(13,5).Range.Text = "18+1=19"
$t.Cell(14,1).Range.Text = "46+39=85"
$t.Cell(14,2).Range.Text = "88-77=11"
$t.Cell(14,3).Range.Text = "29+44=73"
$t.Cell(14,4).Range.Text = "86-19=67"
$t.Cell(14,5).Range.Text = "52+25=77"
$t.Cell(15,1).Range.Text = "42+33=75"
$t.Cell(15,2).Range.Text = "23+69=92"
$t.Cell(15,3).Range.Text = "17+31=48"
$t.Cell(15,4).Range.Text = "36-22=14"
$t.Cell(15,5).Range.Text = "33+46=79"
$t.Cell(16,1).Range.Text = "78-78=0"
$t.Cell(16,2).Range.Text = "56-9=47"
$t.Cell(16,3).Range.Text = "67+13=80"
$t.Cell(16,4).Range.Text = "29-13=16"
$t.Cell(16,5).Range.Text = "71-22=49"
$t.Cell(17,1).Range.Text = "43+30=73"
$t.Cell(17,2).Range.Text = "18+32=50"
$t.Cell(17,3).Range.Text = "41+39=80"
$t.Cell(17,4).Range.Text = "67-32=35"
$t.Cell(17,5).Range.Text = "32+45=77"
$t.Cell(18,1).Range.Text = "11+74=85"
$t.Cell(18,2).Range.Text = "20+49=69"
$t.Cell(18,3).Range.Text = "64-46=18"
$t.Cell(18,4).Range.Text = "36+12=48"
$t.Cell(18,5).Range.Text = "78-37=41"
$t.Cell(19,1).Range.Text = "64-51=13"
$t.Cell(19,2).Range.Text = "37-15=22"
$t.Cell(19,3).Range.Text = "81+16=97"
$t.Cell(19,4).Range.Text = "22+64=86"
$t.Cell(19,5).Range.Text = "60-4=56"
$t.Cell(20,1).Range.Text = "43-32=11"
$t.Cell(20,2).Range.Text = "55+14=69"
$t.Cell(20,3).Range.Text = "48-3=45"
$t.Cell(20,4).Range.Text = "69-35=34"
$t.Cell(20,5).Range.Text = "15+14=29"
